$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9544851183891296
$ws.Range("B1").Value = 1.581608891487122
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.609119653701782
$ws.Range("E1").Value = 1.360632658004761
